$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "26.303.51"
Set-TextValue $ws.Range("E2") "  -0.02%  "
Set-TextValue $ws.Range("D3") "1.690.78"
Set-TextValue $ws.Range("E3") "  +0.65%  "
Set-TextValue $ws.Range("E4") "  +0.15%  "
Set-TextValue $ws.Range("D5") "217.88"
Set-TextValue $ws.Range("E5") "  -0.24%  "
Set-TextValue $ws.Range("D6") "0.5383"
Set-TextValue $ws.Range("E6") "  +2.55%  "
Set-TextValue $ws.Range("E7") "  +0.13%  "
Set-TextValue $ws.Range("E8") "  +1.23%  "
Set-TextValue $ws.Range("D9") "0.06441"
Set-TextValue $ws.Range("E9") "  -0.54%  "
Set-TextValue $ws.Range("D10") "21.64"
Set-TextValue $ws.Range("E10") "  -1.60%  "
Set-TextValue $ws.Range("D11") "0.07669"
Set-TextValue $ws.Range("E11") "  +1.84%  "
Set-TextValue $ws.Range("D12") "1.690.71"
Set-TextValue $ws.Range("E12") "  +0.69%  "
Set-TextValue $ws.Range("D13") "4.535"
Set-TextValue $ws.Range("E13") "  +0.05%  "
Set-TextValue $ws.Range("D14") "0.5789"
Set-TextValue $ws.Range("E14") "  -0.28%  "
Set-TextValue $ws.Range("D15") "0.000008375"
Set-TextValue $ws.Range("D16") "66.98"
Set-TextValue $ws.Range("E16") "  +3.43%  "
Set-TextValue $ws.Range("D17") "26.348.73"
Set-TextValue $ws.Range("E17") "  +0.01%  "
Set-TextValue $ws.Range("D18") "4.908"
Set-TextValue $ws.Range("E18") "  -0.35%  "
Set-TextValue $ws.Range("E19") "  +0.17%  "
Set-TextValue $ws.Range("D20") "10.86"
Set-TextValue $ws.Range("E20") "  -0.09%  "
Set-TextValue $ws.Range("D21") "190.13"
Set-TextValue $ws.Range("E21") "  -0.04%  "
Set-TextValue $ws.Range("D22") "6.265"
Set-TextValue $ws.Range("E22") "  +0.88%  "
Set-TextValue $ws.Range("E23") "  +0.14%  "
Set-TextValue $ws.Range("D24") "149.08"
Set-TextValue $ws.Range("E24") "  +2.34%  "
Set-TextValue $ws.Range("D25") "0.1286"
Set-TextValue $ws.Range("E25") "  +3.05%  "
Set-TextValue $ws.Range("D26") "7.865"
Set-TextValue $ws.Range("E26") "  +0.69%  "
Set-TextValue $ws.Range("D27") "15.85"
Set-TextValue $ws.Range("E27") "  +0.37%  "
Set-TextValue $ws.Range("D28") "0.06264"
Set-TextValue $ws.Range("E28") "  -2.97%  "
Set-TextValue $ws.Range("D29") "1.372"
Set-TextValue $ws.Range("E29") "  +0.82%  "
Set-TextValue $ws.Range("E30") "  +0.06%  "
Set-TextValue $ws.Range("E31") "  -0.07%  "
Set-TextValue $ws.Range("D32") "3.583"
Set-TextValue $ws.Range("E32") "  -0.37%  "
Set-TextValue $ws.Range("D33") "1.674"
Set-TextValue $ws.Range("E33") "  +0.43%  "
Set-TextValue $ws.Range("D34") "1.032"
Set-TextValue $ws.Range("E34") "  +0.13%  "
Set-TextValue $ws.Range("D35") "0.6146"
Set-TextValue $ws.Range("E35") "  -1.66%  "
Set-TextValue $ws.Range("D36") "2.416"
Set-TextValue $ws.Range("E36") "  +0.47%  "
Set-TextValue $ws.Range("D37") "2.768"
Set-TextValue $ws.Range("E37") "  +1.93%  "
Set-TextValue $ws.Range("D38") "0.01653"
Set-TextValue $ws.Range("E38") "  +1.73%  "
Set-TextValue $ws.Range("D39") "1.109.45"
Set-TextValue $ws.Range("E39") "  -0.16%  "
Set-TextValue $ws.Range("E40") "  -5.11%  "
Set-TextValue $ws.Range("D41") "0.8821"
Set-TextValue $ws.Range("E41") "  +0.71%  "
Set-TextValue $ws.Range("D43") "101.38"
Set-TextValue $ws.Range("E43") "  +0.70%  "
Set-TextValue $ws.Range("D44") "1.842.85"
Set-TextValue $ws.Range("E44") "  +0.67%  "
Set-TextValue $ws.Range("D45") "0.00000000113"
Set-TextValue $ws.Range("E45") "  +1.01%  "
Set-TextValue $ws.Range("D46") "57.62"
Set-TextValue $ws.Range("E46") "  +1.14%  "
Set-TextValue $ws.Range("D47") "8.151"
Set-TextValue $ws.Range("E47") "  -0.72%  "
Set-TextValue $ws.Range("D48") "1.001"
Set-TextValue $ws.Range("E48") "  -0.74%  "
Set-TextValue $ws.Range("D49") "0.05282"
Set-TextValue $ws.Range("E49") "  +0.23%  "
Set-TextValue $ws.Range("E50") "  +0.18%  "
Set-TextValue $ws.Range("D51") "6.042"
Set-TextValue $ws.Range("E51") "  -0.73%  "
